$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns (by letter) that vary per price-quote row; A,B,C,E,F,G,H,I,J are
# constant across rows 2-10 for this product/market and are left untouched.
$cols = @("D", "K", "L", "M", "N", "O", "P", "Q", "R", "S", "T")

# Snapshot the "before" values for every row (2-10) and every moving column
# before any writes happen, so the permutation below doesn't clobber data
# it still needs to read.
$snapshot = @{}
for ($r = 2; $r -le 10; $r++) {
    $rowVals = @{}
    foreach ($col in $cols) {
        $rowVals[$col] = $ws.Range("$col$r").Value2
    }
    $snapshot[$r] = $rowVals
}

# The edit re-orders the weekly price quotes: row N's new content is the
# old content that used to live in row `mapping[N]`.
$mapping = @{
    2  = 3
    3  = 4
    4  = 7
    5  = 9
    6  = 10
    7  = 5
    8  = 6
    9  = 8
    10 = 2
}

foreach ($newRow in $mapping.Keys) {
    $oldRow = $mapping[$newRow]
    $src = $snapshot[$oldRow]
    foreach ($col in $cols) {
        $ws.Range("$col$newRow").Value2 = $src[$col]
    }
}
